$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '68.237.37'
$ws.Cells.Item(2, 5).Value = '  -1.08%  '
$ws.Cells.Item(3, 4).Value = '3.898.05'
$ws.Cells.Item(3, 5).Value = '  -0.94%  '
$ws.Cells.Item(4, 5).Value = '  -0.01%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '485.05'
$ws.Cells.Item(5, 4).NumberFormat = "General"
$ws.Cells.Item(5, 5).Value = '  +0.09%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '145.76'
$ws.Cells.Item(6, 4).NumberFormat = "General"
$ws.Cells.Item(6, 5).Value = '  -0.14%  '
$ws.Cells.Item(7, 5).Value = '  -0.17%  '
$ws.Cells.Item(8, 5).Value = '  +0.06%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.740'
$ws.Cells.Item(9, 4).NumberFormat = "General"
$ws.Cells.Item(9, 5).Value = '  +2.33%  '
$ws.Cells.Item(10, 5).Value = '  +7.32%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.0000354'
$ws.Cells.Item(11, 4).NumberFormat = "General"
$ws.Cells.Item(11, 5).Value = '  -0.07%  '
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '42.89'
$ws.Cells.Item(12, 4).NumberFormat = "General"
$ws.Cells.Item(12, 5).Value = '  +0.42%  '
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '10.46'
$ws.Cells.Item(13, 4).NumberFormat = "General"
$ws.Cells.Item(13, 5).Value = '  -0.27%  '
$ws.Cells.Item(14, 4).Value = '4.518.23'
$ws.Cells.Item(14, 5).Value = '  -1.21%  '
$ws.Cells.Item(15, 4).Value = '3.904.65'
$ws.Cells.Item(15, 5).Value = '  -0.92%  '
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '14.17'
$ws.Cells.Item(16, 4).NumberFormat = "General"
$ws.Cells.Item(16, 5).Value = '  -2.79%  '
$ws.Cells.Item(17, 5).Value = '  -0.64%  '
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '19.98'
$ws.Cells.Item(18, 4).NumberFormat = "General"
$ws.Cells.Item(18, 5).Value = '  +1.48%  '
$ws.Cells.Item(19, 5).Value = '  +0.33%  '
$ws.Cells.Item(20, 4).Value = '68.351.34'
$ws.Cells.Item(20, 5).Value = '  -0.98%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '430.43'
$ws.Cells.Item(21, 4).NumberFormat = "General"
$ws.Cells.Item(21, 5).Value = '  -0.87%  '
$ws.Cells.Item(22, 5).Value = '  +6.32%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '12.46'
$ws.Cells.Item(24, 4).NumberFormat = "General"
$ws.Cells.Item(24, 5).Value = '  +22.35%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '89.09'
$ws.Cells.Item(25, 4).NumberFormat = "General"
$ws.Cells.Item(25, 5).Value = '  +1.56%  '
$ws.Cells.Item(26, 5).Value = '  +3.60%  '
$ws.Cells.Item(27, 5).Value = '  -5.42%  '
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '37.26'
$ws.Cells.Item(28, 4).NumberFormat = "General"
$ws.Cells.Item(28, 5).Value = '  -2.70%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '5.69'
$ws.Cells.Item(29, 4).NumberFormat = "General"
$ws.Cells.Item(29, 5).Value = '  -3.74%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '722.39'
$ws.Cells.Item(30, 4).NumberFormat = "General"
$ws.Cells.Item(30, 5).Value = '  +1.36%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '13.40'
$ws.Cells.Item(31, 4).NumberFormat = "General"
$ws.Cells.Item(31, 5).Value = '  +0.77%  '
$ws.Cells.Item(32, 5).Value = '  +0.55%  '
$ws.Cells.Item(33, 5).Value = '  +2.50%  '
$ws.Cells.Item(34, 4).Value = '0.0₃0887'
$ws.Cells.Item(34, 5).Value = '  -3.23%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '61.68'
$ws.Cells.Item(35, 4).NumberFormat = "General"
$ws.Cells.Item(35, 5).Value = '  +4.99%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '6.08'
$ws.Cells.Item(36, 4).NumberFormat = "General"
$ws.Cells.Item(36, 5).Value = '  +7.16%  '
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '40.75'
$ws.Cells.Item(37, 4).NumberFormat = "General"
$ws.Cells.Item(37, 5).Value = '  -1.56%  '
$ws.Cells.Item(38, 5).Value = '  -2.63%  '
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.399'
$ws.Cells.Item(39, 4).NumberFormat = "General"
$ws.Cells.Item(39, 5).Value = '  +17.38%  '
$ws.Cells.Item(40, 5).Value = '  -0.19%  '
$ws.Cells.Item(41, 5).Value = '  +9.43%  '
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '0.0493'
$ws.Cells.Item(42, 4).NumberFormat = "General"
$ws.Cells.Item(42, 5).Value = '  +4.47%  '
$ws.Cells.Item(43, 5).Value = '  +2.97%  '
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '3.02'
$ws.Cells.Item(44, 4).NumberFormat = "General"
$ws.Cells.Item(44, 5).Value = '  -1.34%  '
$ws.Cells.Item(45, 2).Value = 'Stellar'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '0.142'
$ws.Cells.Item(45, 4).NumberFormat = "General"
$ws.Cells.Item(45, 5).Value = '  +0.49%  '
$ws.Cells.Item(46, 2).Value = 'BabyDogeCoin'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Cells.Item(46, 4).Value = '0.0₆0367'
$ws.Cells.Item(46, 5).Value = '  +29.54%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '3.36'
$ws.Cells.Item(47, 4).NumberFormat = "General"
$ws.Cells.Item(47, 5).Value = '  +7.26%  '
$ws.Cells.Item(48, 5).Value = '  -0.03%  '
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '3.36'
$ws.Cells.Item(49, 4).NumberFormat = "General"
$ws.Cells.Item(49, 5).Value = '  -1.29%  '
$ws.Cells.Item(50, 5).Value = '  -2.55%  '
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '143.96'
$ws.Cells.Item(51, 4).NumberFormat = "General"
$ws.Cells.Item(51, 5).Value = '  -2.77%  '
